# Apply "Final FIM March 2023" updates to the K:S (2023 Q1 - 2025 Q1) columns
# for the updated "current" and "difference" rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("K13").Value = -0.0784
$ws.Range("L13").Value = 0.0017
$ws.Range("M13").Value = -0.0419
$ws.Range("N13").Value = -0.0414
$ws.Range("O13").Value = -0.0628
$ws.Range("P13").Value = -0.024
$ws.Range("Q13").Value = 0.0267
$ws.Range("R13").Value = 0.0308
$ws.Range("S13").Value = 0.041

# Row 15
$ws.Range("K15").Value = -0.1448
$ws.Range("L15").Value = 0.0079
$ws.Range("M15").Value = 0.0625
$ws.Range("N15").Value = 0.1251
$ws.Range("O15").Value = 0.4023
$ws.Range("P15").Value = 0.3972
$ws.Range("Q15").Value = 0.4294
$ws.Range("R15").Value = 0.3823
$ws.Range("S15").Value = 0.2794

# Row 31
$ws.Range("K31").Value = -1.4594
$ws.Range("L31").Value = -1.5607
$ws.Range("M31").Value = -0.4575
$ws.Range("N31").Value = -0.4824
$ws.Range("O31").Value = -0.5423
$ws.Range("P31").Value = -0.1539
$ws.Range("Q31").Value = -0.1131
$ws.Range("R31").Value = 0.0286
$ws.Range("S31").Value = 0.0064

# Row 69
$ws.Range("K69").Value = 0.031
$ws.Range("L69").Value = 0.0709
$ws.Range("M69").Value = 0.0956
$ws.Range("N69").Value = 0.0972
$ws.Range("O69").Value = 0.0671
$ws.Range("P69").Value = 0.0294
$ws.Range("Q69").Value = 0.0068
$ws.Range("R69").Value = 0.0068
$ws.Range("S69").Value = 0.0069

# Row 71
$ws.Range("K71").Value = -0.0222
$ws.Range("L71").Value = -0.0331
$ws.Range("M71").Value = -0.0307
$ws.Range("N71").Value = -0.0255
$ws.Range("O71").Value = -0.0203
$ws.Range("P71").Value = -0.0202
$ws.Range("Q71").Value = -0.02
$ws.Range("R71").Value = -0.0209
$ws.Range("S71").Value = -0.0093

# Row 87
$ws.Range("K87").Value = 0.0287
$ws.Range("L87").Value = 0.0454
$ws.Range("M87").Value = 0.0821
$ws.Range("N87").Value = 0.08
$ws.Range("O87").Value = 0.0545
$ws.Range("P87").Value = 0.0236
$ws.Range("Q87").Value = -0.0076
$ws.Range("R87").Value = -0.0143
$ws.Range("S87").Value = -0.0025
